# Tech report cover page update:
#   "Canadian Technical Report of Fisheries and Aquatic Sciences"
#     -> "Canadian Technical Report of Hydrography and Ocean Sciences"
#   "Rapport technique canadien des sciences halieutiques et aquatiques"
#     -> "Rapport technique canadien sur l'hydrographie et les sciences
#         océaniques" (second half set bold, as in the target)
#
# Both target paragraphs live in the small-print legend below the cover
# table; the cover table itself already contains coincidentally similar
# text ("Canadian Technical Report of Hydrography and Ocean Sciences"
# title block with a series-number form field), so every lookup below is
# scoped to the document range that starts *after* the table to avoid
# touching the wrong text.

$d = $word.ActiveDocument

$tbl = $d.Tables(1)
$afterTbl = $d.Range($tbl.Range.End, $d.Content.End)

# ---------------------------------------------------------------------
# 1) "Canadian Technical Report of Fisheries and Aquatic Sciences"
# ---------------------------------------------------------------------
$oldEn = "Canadian Technical Report of Fisheries and Aquatic Sciences"
$prefixEn = "Canadian Technical Report of "
$newSuffixEn = "Hydrography and Ocean Sciences"

$rngEn = $afterTbl.Duplicate
$rngEn.Find.Execute($oldEn, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$suffixEn = $rngEn.Duplicate
$suffixEn.MoveStart(1, $prefixEn.Length)
# Force the replaced text onto its own run (matching the source edit,
# which left the "Canadian Technical Report of " run untouched and
# introduced a new run for the replacement words) by toggling Bold off
# then back on around the text swap.
$suffixEn.Font.Bold = $false
$suffixEn.Text = $newSuffixEn
$suffixEn.Font.Bold = $true

# ---------------------------------------------------------------------
# 2) "Rapport technique canadien des sciences halieutiques et aquatiques"
# ---------------------------------------------------------------------
$oldFr = "Rapport technique canadien des sciences halieutiques et aquatiques"
$prefixFr = "Rapport technique canadien "
$newSuffixFr = "sur l'hydrographie et les sciences océaniques"

$rngFr = $afterTbl.Duplicate
$rngFr.Find.Execute($oldFr, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$suffixFr = $rngFr.Duplicate
$suffixFr.MoveStart(1, $prefixFr.Length)
$suffixFr.Font.Bold = $false
$suffixFr.Text = $newSuffixFr
$suffixFr.Font.Bold = $true
$suffixFr.Font.BoldBi = $true
